$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# "Periodo Mora" column (E) values for the 4 worker rows were re-entered in
# the opposite order, which made Excel rebuild the shared-string table with
# 2402/2406/2407/2408 in ascending order instead of the original
# descending 2408/2407/2406/2402 order.
$ws.Range("E16").Value = "2402"
$ws.Range("E17").Value = "2406"
$ws.Range("E18").Value = "2407"
$ws.Range("E19").Value = "2408"

# "Valor Mora" column (F) amounts for the first and last rows were swapped.
$ws.Range("F16").Value = 29466
$ws.Range("F19").Value = 42000
